$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leading "=" from the G6 message text so it reads as a literal
# string rather than beginning with a formula-like sign.
$ws.Range("G6").Value = "Message A6, E, C6, IF(0=10,TRUE,FALSE), false"

# Move the active selection from G11 (outside the used range) to G7.
$ws.Range("G7").Select()

# Shift the window's vertical position.
$excel.ActiveWindow.Top = $excel.ActiveWindow.Top + 450
